$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text storage (matching the original inlineStr cell type) by temporarily
# applying a text number format before assignment, then restoring the
# default "Normal" style so no stray formatting is left behind.
$numericLookingCells = @('D4', 'D5', 'D6', 'D7', 'D12', 'D14', 'D21', 'D23', 'D24', 'D25', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D38', 'D39', 'D41', 'D42', 'D45', 'D48', 'D49', 'D50')
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated Coin / Link / Price / Volume(1h) values cell by cell.
$ws.Range('D2').Value = '66.075.15'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '3.314.19'
$ws.Range('E3').Value = '  -1.61%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '187.78'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('D6').Value = '559.76'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = '3.306.28'
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').Value = '47.68'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = '8.63'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('E15').Value = '  +3.92%  '
$ws.Range('D16').Value = '3.844.50'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('E17').Value = '  +2.92%  '
$ws.Range('D18').Value = '66.077.52'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D20').Value = '3.291.87'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').Value = '11.12'
$ws.Range('E21').Value = '  -4.85%  '
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '18.15'
$ws.Range('E23').Value = '  +6.30%  '
$ws.Range('D24').Value = '102.74'
$ws.Range('E24').Value = '  +7.31%  '
$ws.Range('D25').Value = '4.96'
$ws.Range('E25').Value = '  -2.59%  '
$ws.Range('E26').Value = '  -3.86%  '
$ws.Range('D27').Value = '6.05'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('D29').Value = '9.56'
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = '8.68'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').Value = '30.28'
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('D32').Value = '4.06'
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('D33').Value = '6.35'
$ws.Range('E33').Value = '  -0.15%  '
$ws.Range('D34').Value = '11.08'
$ws.Range('E34').Value = '  -1.65%  '
$ws.Range('D35').Value = '551.95'
$ws.Range('E35').Value = '  +3.47%  '
$ws.Range('D36').Value = '3.851.90'
$ws.Range('E36').Value = '  +2.57%  '
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').Value = '57.64'
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').Value = '0.0₃0733'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '3.30'
$ws.Range('E41').Value = '  -3.41%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').Value = '33.75'
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('E44').Value = '  -1.23%  '
$ws.Range('D45').Value = '3.23'
$ws.Range('E45').Value = '  -14.08%  '
$ws.Range('E46').Value = '  -5.18%  '
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = '3.24'
$ws.Range('E48').Value = '  +1.62%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').Value = '2.61'
$ws.Range('E49').Value = '  -3.03%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').Value = '0.129'
$ws.Range('E50').Value = '  -1.04%  '
$ws.Range('E51').Value = '  -0.02%  '

# Restore default styling on the cells we temporarily reformatted.
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).Style = "Normal"
}

Write-Output "Updated cryptos list"
